# ============================================================
# Edit script for A5024_Protocol.docx
# ============================================================
$d = $word.ActiveDocument

# ------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark near "Млечни продукти:"
# ------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------
# 2) Remove the "№255" run from "писмо №255 от ..."
#    (delete just the field text; the whole run collapses away)
# ------------------------------------------------------------
$null = $d.Content.Find.Execute("№255", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 2)

# ------------------------------------------------------------
# 3) Reorder the BDS standards list
# ------------------------------------------------------------
$null = $d.Content.Find.Execute("БДС 000-2; БДС 123; БДС 123-2", $false, $false, $false, `
    $false, $false, $true, 1, $false, "БДС 000-2; БДС 123-2; БДС 123", 2)

# ------------------------------------------------------------
# 4) Change the time "18:22" to "21:22" (only the "18" run)
# ------------------------------------------------------------
$null = $d.Content.Find.Execute("18", $false, $false, $false, $false, $false, `
    $true, 1, $false, "21", 1)

# ------------------------------------------------------------
# 5) Replace the "#REMARKSLIST" merge-field paragraph with the
#    literal remarks text, and drop the tiny blank spacer
#    paragraph (font size 6pt / sz=12) right above it.
# ------------------------------------------------------------
$n = $d.Paragraphs.Count
$targetIdx = -1
for ($i = 1; $i -le $n; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "REMARKSLIST") {
        $targetIdx = $i
        break
    }
}

$emptyPara = $d.Paragraphs($targetIdx - 1)
$null = $emptyPara.Range.Delete()

$n2 = $d.Paragraphs.Count
$targetIdx2 = -1
for ($i = 1; $i -le $n2; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "REMARKSLIST") {
        $targetIdx2 = $i
        break
    }
}

$remarksPara = $d.Paragraphs($targetIdx2)
$remarksXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="120"/><w:jc w:val="both"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>
Забележка 1: Lorem Ipsum is simply dummy text of the printing and typesetting industry. Lorem Ipsum has been the industry's standard dummy text ever since
</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>

</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>
Забележка 2: Забележка 4
</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>

</w:t></w:r><w:r><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $remarksPara.Range.InsertXML($remarksXml)

# ------------------------------------------------------------
# 6) Update the cached "Страница X от Y" PAGE field result in
#    the header from 2 to 1 (leave NUMPAGES field untouched)
# ------------------------------------------------------------
$hdr = $d.Sections(1).Headers(1)
$null = $hdr.Range.Find.Execute("2", $false, $false, $false, $false, $false, `
    $true, 1, $false, "1", 1)
